# Add a "Unit" column to the import template.
# The new column is inserted as column I (between "PO Number" -> now J,
# and the existing columns shift right by one), matching the target
# layout: A..H unchanged, I="Unit", J="PO Number", K="Item category",
# L="Item subcategory", M.."Extra Details" shifted to M..W.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (9th column), shifting existing
# columns I.. to the right.
$ws.Columns.Item(9).Insert()

# Set the header text for the newly inserted column.
$ws.Cells.Item(1, 9).Value = "Unit"
